# Update the call-center report figures on the active sheet.
# The underlying numbers (llamada, Real Gestionados, Real Contactados,
# Real Contactos Efectivos, Real Valoraciones Positivas) were recalculated
# "ahora funciona en todos lados" - apply the refreshed values per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = 21;  "G2" = 15;  "I2" = 6;   "K2" = 1;   "M2" = 1
    "D3" = 65;  "E3" = 26;  "G3" = 52;  "I3" = 27;  "K3" = 24
    "D4" = 12;  "E4" = 6;   "G4" = 9
    "D5" = 53;  "E5" = 35;  "F5" = 11;  "G5" = 45;  "I5" = 16;  "K5" = 6
    "D7" = 35;  "E7" = 23;  "F7" = 1;   "G7" = 30;  "I7" = 8;   "K7" = 6;  "M7" = 4
    "D8" = 15;  "E8" = 14;  "G8" = 14
    "D9" = 52;  "E9" = 39;  "F9" = 3;   "G9" = 45;  "I9" = 9;   "K9" = 6
    "D10" = 53; "F10" = 4;  "G10" = 31; "I10" = 6;  "K10" = 2;  "M10" = 2
    "D11" = 12; "F11" = 5;  "G11" = 10; "I11" = 10; "K11" = 6;  "M11" = 5
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
